$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 2

$ws.Range("G9").Select()
